$wb = $excel.ActiveWorkbook

# --- Sheet "Introduction " updates ---
$wsIntro = $wb.Worksheets.Item("Introduction ")

# Ref Version: 1.1 -> 1.3
$wsIntro.Range("D7").Value = 1.3

$lastUpdateDate = Get-Date -Year 2020 -Month 7 -Day 2 -Hour 0 -Minute 0 -Second 0

# Last update (D9): turn the plain text date into a real date value
$wsIntro.Range("D9").Value = $lastUpdateDate
$wsIntro.Range("D9").NumberFormat = "mm-dd-yy"

# New revision-history row (row 14)
$wsIntro.Range("B14").Value = 0.2
$wsIntro.Range("C14").Value = "T.Sharaby"
$wsIntro.Range("E14").Value = $lastUpdateDate
$wsIntro.Range("G14").Value = "Update the status "

# Re-use the same date style for E14 as D9 (so both share one style record)
$wsIntro.Range("D9").Copy()
$wsIntro.Range("E14").PasteSpecial(-4122)

# Select B10:H10 on the Introduction sheet and make it the active sheet/tab
$wsIntro.Range("B10:H10").Select()
$wsIntro.Activate()

# --- Sheet "Cross review points " updates ---
$wsReview = $wb.Worksheets.Item("Cross review points ")

# Status column: mark the first batch of open points as Resolved
$wsReview.Range("H2:H8").Value = "Resolved"

$wsReview.Range("F11").Select()

$wsIntro.Activate()
